$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay as text (avoid numeric auto-conversion)

# D/E price and volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.308.41"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.89"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.52"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3776"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  -4.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.20"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.006"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.283"
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08147"
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.92"
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.684"
$ws.Range("E14").Value = "  -6.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.574"
$ws.Range("E15").Value = "  -4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001272"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.600.79"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.55"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06828"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.62"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.647"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.301.21"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.397"
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.975"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.25"
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.92"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.324"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.26"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.776.83"
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.373"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2559"
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08893"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.397"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7211"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.88"
$ws.Range("E43").Value = "  -4.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.00"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6644"
$ws.Range("E45").Value = "  -4.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.326"
$ws.Range("E46").Value = "  -5.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.981"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08062"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.73"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.181"
$ws.Range("E51").Value = "  -4.17%  "

# Volume-only updates (price unchanged)
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("E36").Value = "  -5.87%  "
$ws.Range("E39").Value = "  -6.86%  "

# Row reorders with updated data (31<->32 swap content, 34<->35 swap content)
$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.467"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.509"
$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9760"
$ws.Range("E34").Value = "  -6.19%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07791"
$ws.Range("E35").Value = "  -3.51%  "
